$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 275.7
$ws.Cells.Item(2, 9).Value = 209.83333
$ws.Cells.Item(2, 10).Value = 374.5
$ws.Cells.Item(2, 11).Value = 209.83333
$ws.Cells.Item(2, 12).Value = 374.5
$ws.Cells.Item(2, 13).Value = -96.83332999999999
$ws.Cells.Item(2, 14).Value = -600.5
$ws.Cells.Item(76, 8).Value = 3780.6365
$ws.Cells.Item(76, 9).Value = 3161.625
$ws.Cells.Item(76, 11).Value = 3161.625
$ws.Cells.Item(76, 13).Value = -2846.625
$ws.Cells.Item(79, 8).Value = 3780.6365
$ws.Cells.Item(79, 9).Value = 3161.625
$ws.Cells.Item(79, 11).Value = 3161.625
$ws.Cells.Item(79, 13).Value = -2069.625
$ws.Cells.Item(80, 8).Value = 1348.7333
$ws.Cells.Item(80, 9).Value = 1906.75
$ws.Cells.Item(80, 11).Value = 5720.25
$ws.Cells.Item(80, 13).Value = -4722.25
$ws.Cells.Item(83, 8).Value = 1348.7333
$ws.Cells.Item(83, 9).Value = 1906.75
$ws.Cells.Item(83, 11).Value = 17160.75
$ws.Cells.Item(83, 13).Value = -12168.75
$ws.Cells.Item(132, 8).Value = 950.8163500000001
$ws.Cells.Item(132, 9).Value = 942.24445
$ws.Cells.Item(132, 10).Value = 1047.25
$ws.Cells.Item(132, 11).Value = 2826.73335
$ws.Cells.Item(132, 12).Value = 3141.75
$ws.Cells.Item(132, 13).Value = -296.73335
$ws.Cells.Item(132, 14).Value = -8201.75
$ws.Cells.Item(137, 8).Value = 56687.223
$ws.Cells.Item(137, 9).Value = 771.2222
$ws.Cells.Item(137, 11).Value = 2313.6666
$ws.Cells.Item(137, 13).Value = 236.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2485.2688
$ws.Cells.Item(32, 9).Value = 1966.7108
$ws.Cells.Item(32, 11).Value = 1966.7108
$ws.Cells.Item(32, 13).Value = -1679.7108
$ws.Cells.Item(35, 8).Value = 3778.6667
$ws.Cells.Item(35, 9).Value = 3778.6667
$ws.Cells.Item(35, 11).Value = 3778.6667
$ws.Cells.Item(35, 13).Value = -3372.6667
$ws.Cells.Item(63, 8).Value = 3179.3333
$ws.Cells.Item(63, 9).Value = 3515.2
$ws.Cells.Item(63, 10).Value = 1500
$ws.Cells.Item(63, 11).Value = 3515.2
$ws.Cells.Item(63, 12).Value = 1500
$ws.Cells.Item(63, 13).Value = -2829.2
$ws.Cells.Item(63, 14).Value = -2872
$ws.Cells.Item(66, 8).Value = 3179.3333
$ws.Cells.Item(66, 9).Value = 3515.2
$ws.Cells.Item(66, 10).Value = 1500
$ws.Cells.Item(66, 11).Value = 17576
$ws.Cells.Item(66, 12).Value = 7500
$ws.Cells.Item(66, 13).Value = -14144
$ws.Cells.Item(66, 14).Value = -14364
$ws.Cells.Item(74, 8).Value = 875.2
$ws.Cells.Item(74, 9).Value = 804.65
$ws.Cells.Item(74, 11).Value = 804.65
$ws.Cells.Item(74, 13).Value = 69.35000000000002
$ws.Cells.Item(77, 8).Value = 875.2
$ws.Cells.Item(77, 9).Value = 804.65
$ws.Cells.Item(77, 11).Value = 4023.25
$ws.Cells.Item(77, 13).Value = 344.75
$ws.Cells.Item(122, 8).Value = 1674.2632
$ws.Cells.Item(122, 9).Value = 1800.7333
$ws.Cells.Item(122, 11).Value = 5402.199900000001
$ws.Cells.Item(122, 13).Value = -2952.199900000001
$ws.Cells.Item(132, 8).Value = 2355.975
$ws.Cells.Item(132, 9).Value = 2018.7142
$ws.Cells.Item(132, 10).Value = 2728.7368
$ws.Cells.Item(132, 11).Value = 6056.142599999999
$ws.Cells.Item(132, 12).Value = 8186.2104
$ws.Cells.Item(132, 13).Value = -3526.142599999999
$ws.Cells.Item(132, 14).Value = -13246.2104

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 83.333336
$ws.Cells.Item(22, 9).Value = 112.5
$ws.Cells.Item(22, 11).Value = 112.5
$ws.Cells.Item(22, 13).Value = 60.5
$ws.Cells.Item(105, 8).Value = 2669.8572
$ws.Cells.Item(105, 9).Value = 2448.25
$ws.Cells.Item(105, 11).Value = 2448.25
$ws.Cells.Item(105, 13).Value = -701.25
$ws.Cells.Item(134, 8).Value = 8750.959999999999
$ws.Cells.Item(134, 9).Value = 9688.75
$ws.Cells.Item(134, 11).Value = 29066.25
$ws.Cells.Item(134, 13).Value = -26531.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1981.3043
$ws.Cells.Item(31, 9).Value = 1087.5714
$ws.Cells.Item(31, 10).Value = 2732.04
$ws.Cells.Item(31, 11).Value = 1087.5714
$ws.Cells.Item(31, 12).Value = 2732.04
$ws.Cells.Item(31, 13).Value = -792.5714
$ws.Cells.Item(31, 14).Value = -3322.04
$ws.Cells.Item(34, 8).Value = 1981.3043
$ws.Cells.Item(34, 9).Value = 1087.5714
$ws.Cells.Item(34, 10).Value = 2732.04
$ws.Cells.Item(34, 11).Value = 1087.5714
$ws.Cells.Item(34, 12).Value = 2732.04
$ws.Cells.Item(34, 13).Value = -885.5714
$ws.Cells.Item(34, 14).Value = -3136.04
$ws.Cells.Item(107, 8).Value = 372.875
$ws.Cells.Item(107, 9).Value = 340.42856
$ws.Cells.Item(107, 11).Value = 340.42856
$ws.Cells.Item(107, 13).Value = 1579.57144
$ws.Cells.Item(122, 8).Value = 1483.5
$ws.Cells.Item(122, 9).Value = 1006.1111
$ws.Cells.Item(122, 10).Value = 2342.8
$ws.Cells.Item(122, 11).Value = 3018.3333
$ws.Cells.Item(122, 12).Value = 7028.400000000001
$ws.Cells.Item(122, 13).Value = -568.3332999999998
$ws.Cells.Item(122, 14).Value = -11928.4
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 5385.4614
$ws.Cells.Item(132, 9).Value = 4249.75
$ws.Cells.Item(132, 11).Value = 12749.25
$ws.Cells.Item(132, 13).Value = -10219.25
$ws.Range("N124").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 549.5
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(2, 11).Value = 600
$ws.Cells.Item(2, 13).Value = -487
$ws.Cells.Item(7, 8).Value = 399.46155
$ws.Cells.Item(7, 9).Value = 212.875
$ws.Cells.Item(7, 11).Value = 638.625
$ws.Cells.Item(7, 13).Value = -526.625
$ws.Cells.Item(32, 8).Value = 1686.1666
$ws.Cells.Item(32, 10).Value = 1686.1666
$ws.Cells.Item(32, 12).Value = 5058.4998
$ws.Cells.Item(32, 14).Value = -5624.4998
$ws.Cells.Item(39, 8).Value = 2505.3125
$ws.Cells.Item(39, 10).Value = 2922.1538
$ws.Cells.Item(39, 12).Value = 8766.4614
$ws.Cells.Item(39, 14).Value = -9354.4614
$ws.Cells.Item(57, 8).Value = 3500
$ws.Cells.Item(57, 10).Value = 3500
$ws.Cells.Item(57, 12).Value = 10500
$ws.Cells.Item(57, 14).Value = -11618
$ws.Cells.Item(61, 8).Value = 227.75
$ws.Cells.Item(61, 10).Value = 237
$ws.Cells.Item(61, 12).Value = 711
$ws.Cells.Item(61, 14).Value = -1141
$ws.Cells.Item(114, 8).Value = 11906824
$ws.Cells.Item(114, 10).Value = 17859908
$ws.Cells.Item(114, 12).Value = 53579724
$ws.Cells.Item(114, 14).Value = -53586232
$ws.Cells.Item(137, 8).Value = 2821.6562
$ws.Cells.Item(137, 9).Value = 864.8461
$ws.Cells.Item(137, 10).Value = 4160.5264
$ws.Cells.Item(137, 11).Value = 2594.5383
$ws.Cells.Item(137, 12).Value = 12481.5792
$ws.Cells.Item(137, 13).Value = 2505.4617
$ws.Cells.Item(137, 14).Value = -22681.5792
$ws.Cells.Item(138, 8).Value = 1358.375
$ws.Cells.Item(138, 9).Value = 1159.7142
$ws.Cells.Item(138, 10).Value = 2749
$ws.Cells.Item(138, 11).Value = 3479.1426
$ws.Cells.Item(138, 12).Value = 8247
$ws.Cells.Item(138, 13).Value = 1660.8574
$ws.Cells.Item(138, 14).Value = -18527

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3393
$ws.Cells.Item(122, 10).Value = 3592.125
$ws.Cells.Item(122, 12).Value = 10776.375
$ws.Cells.Item(122, 14).Value = -15676.375
$ws.Cells.Item(132, 8).Value = 787186.1
$ws.Cells.Item(132, 9).Value = 1242257.1
$ws.Cells.Item(132, 10).Value = 3452.7778
$ws.Cells.Item(132, 11).Value = 3726771.3
$ws.Cells.Item(132, 12).Value = 10358.3334
$ws.Cells.Item(132, 13).Value = -3724241.3
$ws.Cells.Item(132, 14).Value = -15418.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2074.75
$ws.Cells.Item(22, 9).Value = 1797
$ws.Cells.Item(22, 10).Value = 2167.3333
$ws.Cells.Item(22, 11).Value = 1797
$ws.Cells.Item(22, 12).Value = 2167.3333
$ws.Cells.Item(22, 13).Value = -1502
$ws.Cells.Item(22, 14).Value = -2757.3333
$ws.Cells.Item(27, 8).Value = 2074.75
$ws.Cells.Item(27, 9).Value = 1797
$ws.Cells.Item(27, 10).Value = 2167.3333
$ws.Cells.Item(27, 11).Value = 1797
$ws.Cells.Item(27, 12).Value = 2167.3333
$ws.Cells.Item(27, 13).Value = -1690
$ws.Cells.Item(27, 14).Value = -2381.3333
$ws.Cells.Item(132, 8).Value = 2968.4138
$ws.Cells.Item(132, 9).Value = 1015.8421
$ws.Cells.Item(132, 10).Value = 6678.3
$ws.Cells.Item(132, 11).Value = 3047.5263
$ws.Cells.Item(132, 12).Value = 20034.9
$ws.Cells.Item(132, 13).Value = -517.5263
$ws.Cells.Item(132, 14).Value = -25094.9
$ws.Cells.Item(133, 8).Value = 89000
$ws.Cells.Item(133, 10).Value = 89000
$ws.Cells.Item(133, 12).Value = 89000
$ws.Cells.Item(133, 14).Value = -94060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2474.05
$ws.Cells.Item(81, 9).Value = 2373.8125
$ws.Cells.Item(81, 11).Value = 4747.625
$ws.Cells.Item(81, 13).Value = -3686.625
$ws.Cells.Item(84, 8).Value = 2474.05
$ws.Cells.Item(84, 9).Value = 2373.8125
$ws.Cells.Item(84, 11).Value = 23738.125
$ws.Cells.Item(84, 13).Value = -18434.125
$ws.Cells.Item(132, 8).Value = 1842.9032
$ws.Cells.Item(132, 9).Value = 817.8
$ws.Cells.Item(132, 10).Value = 2331.0476
$ws.Cells.Item(132, 11).Value = 2453.4
$ws.Cells.Item(132, 12).Value = 6993.1428
$ws.Cells.Item(132, 13).Value = 76.60000000000036
$ws.Cells.Item(132, 14).Value = -12053.1428
$ws.Cells.Item(136, 8).Value = 12348117
$ws.Cells.Item(136, 9).Value = 32682468
$ws.Cells.Item(136, 10).Value = 2261.8572
$ws.Cells.Item(136, 11).Value = 98047404
$ws.Cells.Item(136, 12).Value = 6785.571599999999
$ws.Cells.Item(136, 13).Value = -98044854
$ws.Cells.Item(136, 14).Value = -11885.5716
